$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 72.266001
$ws.Range("H2").Value = 216.798003
$ws.Range("I2").Value = 0.2949652269937106
$ws.Range("J2").Value = 0.2949652269937106
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 146.491518705782
$ws.Range("R2").Value = 1318.423668352038
$ws.Range("S2").Value = 0.001945674574153811
$ws.Range("T2").Value = 0.001945674574153811

$ws.Range("G3").Value = 72.266001
$ws.Range("H3").Value = 216.798003
$ws.Range("I3").Value = 0.2949652269937106
$ws.Range("J3").Value = 0.2949652269937106
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 18532.20442566297
$ws.Range("R3").Value = 166789.8398309667
$ws.Range("S3").Value = 0.246141478172893
$ws.Range("T3").Value = 0.2461414781728929

$ws.Range("G4").Value = 72.266001
$ws.Range("H4").Value = 216.798003
$ws.Range("I4").Value = 0.2949652269937106
$ws.Range("J4").Value = 0.2949652269937106
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 3529.4906875076
$ws.Range("R4").Value = 31765.41618756839
$ws.Range("S4").Value = 0.04687807424666383
$ws.Range("T4").Value = 0.04687807424666382

$ws.Range("H5").Value = 410.023338
$ws.Range("I5").Value = 0.5578585839920717
$ws.Range("J5").Value = 0.5578585839920718
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 277.0548651614386
$ws.Range("R5").Value = 2493.493786452948
$ws.Range("S5").Value = 0.003679793967273187
$ws.Range("T5").Value = 0.003679793967273187

$ws.Range("H6").Value = 410.023338
$ws.Range("I6").Value = 0.5578585839920717
$ws.Range("J6").Value = 0.5578585839920718
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("Q6").Value = 35049.38336128817
$ws.Range("R6").Value = 315444.4502515935
$ws.Range("S6").Value = 0.4655197423599133
$ws.Range("T6").Value = 0.4655197423599133

$ws.Range("H7").Value = 410.023338
$ws.Range("I7").Value = 0.5578585839920717
$ws.Range("J7").Value = 0.5578585839920718
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 6675.216252484488
$ws.Range("R7").Value = 60076.94627236039
$ws.Range("S7").Value = 0.08865904766488526
$ws.Range("T7").Value = 0.08865904766488525

$ws.Range("G8").Value = 36.057927
$ws.Range("H8").Value = 108.173781
$ws.Range("I8").Value = 0.1471761890142177
$ws.Range("J8").Value = 0.1471761890142177
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 73.09357670991399
$ws.Range("R8").Value = 657.8421903892259
$ws.Range("S8").Value = 0.000970816023991617
$ws.Range("T8").Value = 0.0009708160239916168

$ws.Range("G9").Value = 36.057927
$ws.Range("H9").Value = 108.173781
$ws.Range("I9").Value = 0.1471761890142177
$ws.Range("J9").Value = 0.1471761890142177
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 9246.850041275044
$ws.Range("R9").Value = 83221.65037147538
$ws.Range("S9").Value = 0.1228150351315312
$ws.Range("T9").Value = 0.1228150351315312

$ws.Range("G10").Value = 36.057927
$ws.Range("H10").Value = 108.173781
$ws.Range("I10").Value = 0.1471761890142177
$ws.Range("J10").Value = 0.1471761890142177
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 1761.078734069273
$ws.Range("R10").Value = 15849.70860662346
$ws.Range("S10").Value = 0.02339033785869491
$ws.Range("T10").Value = 0.0233903378586949
